$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set BF column (Date) cells to Text format to prevent auto date-conversion
$ws.Range("BF2:BF31").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = 82
$ws.Range("F2").Value = 34
$ws.Range("G2").Value = 0.585
$ws.Range("I2").Value = 38.6
$ws.Range("J2").Value = 84.40000000000001
$ws.Range("L2").Value = 9.9
$ws.Range("N2").Value = 0.35
$ws.Range("O2").Value = 15.6
$ws.Range("U2").Value = 25.6
$ws.Range("X2").Value = 5.9
$ws.Range("AB2").Value = 102.8
$ws.Range("AC2").Value = 3.6
$ws.Range("AI2").Value = 9
$ws.Range("AO2").Value = 29
$ws.Range("AQ2").Value = 8
$ws.Range("AS2").Value = 13
$ws.Range("AV2").Value = 22
$ws.Range("BF2").Value = "2016-06-12"

# Row 3
$ws.Range("D3").Value = 82
$ws.Range("E3").Value = 48
$ws.Range("G3").Value = 0.585
$ws.Range("J3").Value = 89.2
$ws.Range("L3").Value = 8.699999999999999
$ws.Range("O3").Value = 18.5
$ws.Range("Q3").Value = 0.788
$ws.Range("S3").Value = 33.3
$ws.Range("T3").Value = 44.9
$ws.Range("U3").Value = 24.2
$ws.Range("AA3").Value = 21
$ws.Range("AB3").Value = 105.7
$ws.Range("AC3").Value = 3.2
$ws.Range("AE3").Value = 7
$ws.Range("AF3").Value = 7
$ws.Range("AG3").Value = 7
$ws.Range("AH3").Value = 23
$ws.Range("AK3").Value = 24
$ws.Range("AN3").Value = 28
$ws.Range("AO3").Value = 8
$ws.Range("AQ3").Value = 6
$ws.Range("AR3").Value = 3
$ws.Range("AV3").Value = 12
$ws.Range("AX3").Value = 22
$ws.Range("AY3").Value = 24
$ws.Range("BF3").Value = "2016-06-12"

# Row 4
$ws.Range("D4").Value = 82
$ws.Range("F4").Value = 61
$ws.Range("G4").Value = 0.256
$ws.Range("J4").Value = 84.40000000000001
$ws.Range("K4").Value = 0.453
$ws.Range("L4").Value = 6.5
$ws.Range("M4").Value = 18.4
$ws.Range("N4").Value = 0.352
$ws.Range("O4").Value = 15.7
$ws.Range("P4").Value = 20.7
$ws.Range("Q4").Value = 0.757
$ws.Range("Y4").Value = 5.2
$ws.Range("AA4").Value = 18.4
$ws.Range("AC4").Value = -7.4
$ws.Range("AO4").Value = 27
$ws.Range("BF4").Value = "2016-06-12"

# Row 5
$ws.Range("D5").Value = 82
$ws.Range("F5").Value = 34
$ws.Range("G5").Value = 0.585
$ws.Range("I5").Value = 37
$ws.Range("L5").Value = 10.6
$ws.Range("M5").Value = 29.4
$ws.Range("N5").Value = 0.362
$ws.Range("P5").Value = 23.7
$ws.Range("Q5").Value = 0.79
$ws.Range("R5").Value = 9
$ws.Range("T5").Value = 43.9
$ws.Range("U5").Value = 21.7
$ws.Range("V5").Value = 12.5
$ws.Range("X5").Value = 5.3
$ws.Range("Y5").Value = 5.5
$ws.Range("Z5").Value = 18.1
$ws.Range("AA5").Value = 20.4
$ws.Range("AB5").Value = 103.4
$ws.Range("AC5").Value = 2.7
$ws.Range("AH5").Value = 9
$ws.Range("AK5").Value = 27
$ws.Range("AN5").Value = 8
$ws.Range("AT5").Value = 14
$ws.Range("AY5").Value = 21
$ws.Range("BB5").Value = 11
$ws.Range("BF5").Value = "2016-06-12"

# Row 6
$ws.Range("D6").Value = 82
$ws.Range("E6").Value = 42
$ws.Range("G6").Value = 0.512
$ws.Range("H6").Value = 48.5
$ws.Range("J6").Value = 87.40000000000001
$ws.Range("M6").Value = 21.4
$ws.Range("N6").Value = 0.371
$ws.Range("Q6").Value = 0.787
$ws.Range("R6").Value = 11.1
$ws.Range("T6").Value = 46.3
$ws.Range("U6").Value = 22.8
$ws.Range("Z6").Value = 18.8
$ws.Range("AB6").Value = 101.6
$ws.Range("AC6").Value = -1.5
$ws.Range("AD6").Value = 1
$ws.Range("AE6").Value = 14
$ws.Range("AF6").Value = 14
$ws.Range("AG6").Value = 14
$ws.Range("AI6").Value = 11
$ws.Range("AK6").Value = 22
$ws.Range("AN6").Value = 3
$ws.Range("AQ6").Value = 7
$ws.Range("AR6").Value = 9
$ws.Range("AT6").Value = 3
$ws.Range("AU6").Value = 12
$ws.Range("BC6").Value = 18
$ws.Range("BF6").Value = "2016-06-12"

# Row 7
$ws.Range("D7").Value = 82
$ws.Range("E7").Value = 57
$ws.Range("F7").Value = 25
$ws.Range("G7").Value = 0.695
$ws.Range("J7").Value = 84
$ws.Range("K7").Value = 0.46
$ws.Range("O7").Value = 16.3
$ws.Range("P7").Value = 21.7
$ws.Range("Q7").Value = 0.748
$ws.Range("S7").Value = 33.9
$ws.Range("T7").Value = 44.5
$ws.Range("U7").Value = 22.7
$ws.Range("AB7").Value = 104.3
$ws.Range("AC7").Value = 6
$ws.Range("AD7").Value = 1
$ws.Range("AI7").Value = 8
$ws.Range("AJ7").Value = 19
$ws.Range("AL7").Value = 2
$ws.Range("AN7").Value = 7
$ws.Range("AP7").Value = 24
$ws.Range("AS7").Value = 11
$ws.Range("AU7").Value = 13
$ws.Range("BB7").Value = 8
$ws.Range("BF7").Value = "2016-06-12"

# Row 8
$ws.Range("D8").Value = 82
$ws.Range("E8").Value = 42
$ws.Range("G8").Value = 0.512
$ws.Range("I8").Value = 37.4
$ws.Range("M8").Value = 28.6
$ws.Range("N8").Value = 0.344
$ws.Range("P8").Value = 22.3
$ws.Range("Q8").Value = 0.794
$ws.Range("S8").Value = 33.9
$ws.Range("T8").Value = 43.1
$ws.Range("Z8").Value = 19.5
$ws.Range("AA8").Value = 21.4
$ws.Range("AB8").Value = 102.3
$ws.Range("AC8").Value = -0.3
$ws.Range("AE8").Value = 14
$ws.Range("AF8").Value = 14
$ws.Range("AG8").Value = 14
$ws.Range("AN8").Value = 23
$ws.Range("AP8").Value = 21
$ws.Range("AS8").Value = 9
$ws.Range("AT8").Value = 19
$ws.Range("BB8").Value = 16
$ws.Range("BF8").Value = "2016-06-12"

# Row 9
$ws.Range("D9").Value = 82
$ws.Range("E9").Value = 33
$ws.Range("G9").Value = 0.402
$ws.Range("J9").Value = 85.40000000000001
$ws.Range("K9").Value = 0.442
$ws.Range("L9").Value = 8
$ws.Range("N9").Value = 0.338
$ws.Range("Q9").Value = 0.766
$ws.Range("U9").Value = 22.7
$ws.Range("V9").Value = 14.7
$ws.Range("AC9").Value = -3.1
$ws.Range("AE9").Value = 21
$ws.Range("AF9").Value = 21
$ws.Range("AG9").Value = 21
$ws.Range("AK9").Value = 21
$ws.Range("AO9").Value = 10
$ws.Range("AR9").Value = 6
$ws.Range("AS9").Value = 19
$ws.Range("AU9").Value = 14
$ws.Range("AV9").Value = 18
$ws.Range("AX9").Value = 16
$ws.Range("BF9").Value = "2016-06-12"

# Row 10
$ws.Range("D10").Value = 82
$ws.Range("E10").Value = 44
$ws.Range("G10").Value = 0.537
$ws.Range("I10").Value = 37.9
$ws.Range("J10").Value = 86.40000000000001
$ws.Range("K10").Value = 0.439
$ws.Range("M10").Value = 26.2
$ws.Range("N10").Value = 0.345
$ws.Range("O10").Value = 17.1
$ws.Range("P10").Value = 25.5
$ws.Range("Q10").Value = 0.668
$ws.Range("R10").Value = 12.5
$ws.Range("S10").Value = 33.9
$ws.Range("T10").Value = 46.3
$ws.Range("W10").Value = 7
$ws.Range("AA10").Value = 21.6
$ws.Range("AC10").Value = 0.6
$ws.Range("AD10").Value = 1
$ws.Range("AE10").Value = 12
$ws.Range("AG10").Value = 12
$ws.Range("AK10").Value = 25
$ws.Range("AL10").Value = 10
$ws.Range("AN10").Value = 22
$ws.Range("AP10").Value = 5
$ws.Range("AS10").Value = 11
$ws.Range("AT10").Value = 2
$ws.Range("BA10").Value = 5
$ws.Range("BC10").Value = 14
$ws.Range("BF10").Value = "2016-06-12"

# Row 11
$ws.Range("D11").Value = 82
$ws.Range("E11").Value = 73
$ws.Range("G11").Value = 0.89
$ws.Range("I11").Value = 42.5
$ws.Range("J11").Value = 87.3
$ws.Range("K11").Value = 0.487
$ws.Range("L11").Value = 13.1
$ws.Range("N11").Value = 0.416
$ws.Range("O11").Value = 16.7
$ws.Range("Q11").Value = 0.763
$ws.Range("R11").Value = 10
$ws.Range("S11").Value = 36.2
$ws.Range("T11").Value = 46.2
$ws.Range("V11").Value = 15.2
$ws.Range("Z11").Value = 20.7
$ws.Range("AC11").Value = 10.8
$ws.Range("AH11").Value = 5
$ws.Range("AP11").Value = 23
$ws.Range("AR11").Value = 21
$ws.Range("AZ11").Value = 19
$ws.Range("BA11").Value = 19
$ws.Range("BC11").Value = 1
$ws.Range("BF11").Value = "2016-06-12"

# Row 12
$ws.Range("D12").Value = 82
$ws.Range("F12").Value = 41
$ws.Range("G12").Value = 0.5
$ws.Range("I12").Value = 37.7
$ws.Range("K12").Value = 0.452
$ws.Range("M12").Value = 30.9
$ws.Range("N12").Value = 0.347
$ws.Range("Q12").Value = 0.694
$ws.Range("S12").Value = 31.7
$ws.Range("T12").Value = 43.1
$ws.Range("U12").Value = 22.2
$ws.Range("Z12").Value = 21.8
$ws.Range("AB12").Value = 106.5
$ws.Range("AC12").Value = 0.2
$ws.Range("AE12").Value = 17
$ws.Range("AF12").Value = 17
$ws.Range("AG12").Value = 17
$ws.Range("AI12").Value = 20
$ws.Range("AL12").Value = 3
$ws.Range("AT12").Value = 20
$ws.Range("BB12").Value = 4
$ws.Range("BC12").Value = 15
$ws.Range("BF12").Value = "2016-06-12"

# Row 13
$ws.Range("D13").Value = 82
$ws.Range("F13").Value = 37
$ws.Range("G13").Value = 0.549
$ws.Range("I13").Value = 38.3
$ws.Range("K13").Value = 0.45
$ws.Range("N13").Value = 0.351
$ws.Range("O13").Value = 17.4
$ws.Range("P13").Value = 22.8
$ws.Range("Q13").Value = 0.764
$ws.Range("R13").Value = 10.3
$ws.Range("T13").Value = 44.2
$ws.Range("Z13").Value = 20
$ws.Range("AC13").Value = 1.7
$ws.Range("AH13").Value = 5
$ws.Range("AK13").Value = 17
$ws.Range("AN13").Value = 14
$ws.Range("AO13").Value = 14
$ws.Range("AR13").Value = 17
$ws.Range("AS13").Value = 10
$ws.Range("AT13").Value = 12
$ws.Range("AX13").Value = 17
$ws.Range("AY13").Value = 9
$ws.Range("BF13").Value = "2016-06-12"

# Row 14
$ws.Range("D14").Value = 82
$ws.Range("E14").Value = 53
$ws.Range("G14").Value = 0.646
$ws.Range("J14").Value = 82.40000000000001
$ws.Range("K14").Value = 0.465
$ws.Range("L14").Value = 9.699999999999999
$ws.Range("M14").Value = 26.7
$ws.Range("O14").Value = 18.2
$ws.Range("Q14").Value = 0.6919999999999999
$ws.Range("S14").Value = 33.3
$ws.Range("U14").Value = 22.8
$ws.Range("V14").Value = 13
$ws.Range("W14").Value = 8.6
$ws.Range("AA14").Value = 22.5
$ws.Range("AB14").Value = 104.5
$ws.Range("AC14").Value = 4.3
$ws.Range("AK14").Value = 6
$ws.Range("AQ14").Value = 29
$ws.Range("AU14").Value = 11
$ws.Range("BB14").Value = 7
$ws.Range("BF14").Value = "2016-06-12"

# Row 15
$ws.Range("D15").Value = 82
$ws.Range("F15").Value = 65
$ws.Range("G15").Value = 0.207
$ws.Range("M15").Value = 24.6
$ws.Range("N15").Value = 0.317
$ws.Range("O15").Value = 19.3
$ws.Range("P15").Value = 24.7
$ws.Range("Q15").Value = 0.781
$ws.Range("S15").Value = 32.3
$ws.Range("T15").Value = 43
$ws.Range("W15").Value = 7.2
$ws.Range("Z15").Value = 20.3
$ws.Range("AA15").Value = 19.2
$ws.Range("AB15").Value = 97.3
$ws.Range("AC15").Value = -9.6
$ws.Range("AL15").Value = 22
$ws.Range("AP15").Value = 8
$ws.Range("AQ15").Value = 11
$ws.Range("AW15").Value = 23
$ws.Range("BA15").Value = 24
$ws.Range("BF15").Value = "2016-06-12"

# Row 16
$ws.Range("D16").Value = 82
$ws.Range("F16").Value = 40
$ws.Range("G16").Value = 0.512
$ws.Range("I16").Value = 36.8
$ws.Range("J16").Value = 83.59999999999999
$ws.Range("K16").Value = 0.44
$ws.Range("L16").Value = 6.1
$ws.Range("M16").Value = 18.5
$ws.Range("N16").Value = 0.331
$ws.Range("O16").Value = 19.3
$ws.Range("U16").Value = 20.7
$ws.Range("Y16").Value = 5.7
$ws.Range("Z16").Value = 21.7
$ws.Range("AB16").Value = 99.09999999999999
$ws.Range("AC16").Value = -2.2
$ws.Range("AE16").Value = 14
$ws.Range("AI16").Value = 26
$ws.Range("AK16").Value = 23
$ws.Range("AN16").Value = 29
$ws.Range("AP16").Value = 9
$ws.Range("AQ16").Value = 9
$ws.Range("AR16").Value = 8
$ws.Range("AT16").Value = 28
$ws.Range("BF16").Value = "2016-06-12"

# Row 17
$ws.Range("D17").Value = 82
$ws.Range("E17").Value = 48
$ws.Range("G17").Value = 0.585
$ws.Range("J17").Value = 81.7
$ws.Range("L17").Value = 6.1
$ws.Range("N17").Value = 0.336
$ws.Range("P17").Value = 23
$ws.Range("Q17").Value = 0.744
$ws.Range("R17").Value = 9.800000000000001
$ws.Range("S17").Value = 34.3
$ws.Range("T17").Value = 44.1
$ws.Range("W17").Value = 6.7
$ws.Range("Z17").Value = 18.3
$ws.Range("AA17").Value = 19.6
$ws.Range("AC17").Value = 1.6
$ws.Range("AE17").Value = 7
$ws.Range("AF17").Value = 7
$ws.Range("AG17").Value = 7
$ws.Range("AI17").Value = 13
$ws.Range("AN17").Value = 27
$ws.Range("AP17").Value = 14
$ws.Range("AQ17").Value = 23
$ws.Range("AR17").Value = 22
$ws.Range("AV17").Value = 14
$ws.Range("BF17").Value = "2016-06-12"

# Row 18
$ws.Range("D18").Value = 82
$ws.Range("F18").Value = 49
$ws.Range("G18").Value = 0.402
$ws.Range("I18").Value = 38.4
$ws.Range("J18").Value = 82.2
$ws.Range("K18").Value = 0.467
$ws.Range("L18").Value = 5.4
$ws.Range("M18").Value = 15.6
$ws.Range("N18").Value = 0.345
$ws.Range("O18").Value = 17
$ws.Range("P18").Value = 22.7
$ws.Range("Q18").Value = 0.747
$ws.Range("R18").Value = 10.5
$ws.Range("T18").Value = 41.7
$ws.Range("U18").Value = 23.1
$ws.Range("Z18").Value = 20.7
$ws.Range("AA18").Value = 19.5
$ws.Range("AB18").Value = 99
$ws.Range("AC18").Value = -4.2
$ws.Range("AI18").Value = 14
$ws.Range("AN18").Value = 21
$ws.Range("AT18").Value = 27
$ws.Range("AW18").Value = 13
$ws.Range("AZ18").Value = 17
$ws.Range("BA18").Value = 21
$ws.Range("BF18").Value = "2016-06-12"

# Row 19
$ws.Range("D19").Value = 82
$ws.Range("E19").Value = 29
$ws.Range("G19").Value = 0.354
$ws.Range("L19").Value = 5.5
$ws.Range("N19").Value = 0.338
$ws.Range("O19").Value = 21.4
$ws.Range("P19").Value = 27
$ws.Range("Q19").Value = 0.792
$ws.Range("R19").Value = 10
$ws.Range("V19").Value = 15
$ws.Range("W19").Value = 8
$ws.Range("Z19").Value = 20.7
$ws.Range("AA19").Value = 21.6
$ws.Range("AB19").Value = 102.4
$ws.Range("AC19").Value = -3.5
$ws.Range("AH19").Value = 5
$ws.Range("AI19").Value = 19
$ws.Range("AK19").Value = 7
$ws.Range("AN19").Value = 25
$ws.Range("AZ19").Value = 18
$ws.Range("BC19").Value = 24
$ws.Range("BF19").Value = "2016-06-12"

# Row 20
$ws.Range("D20").Value = 82
$ws.Range("F20").Value = 52
$ws.Range("G20").Value = 0.366
$ws.Range("H20").Value = 48.2
$ws.Range("J20").Value = 85.90000000000001
$ws.Range("K20").Value = 0.448
$ws.Range("N20").Value = 0.36
$ws.Range("Q20").Value = 0.776
$ws.Range("R20").Value = 9.5
$ws.Range("S20").Value = 33.1
$ws.Range("T20").Value = 42.6
$ws.Range("AA20").Value = 19.9
$ws.Range("AB20").Value = 102.7
$ws.Range("AC20").Value = -3.8
$ws.Range("AD20").Value = 1
$ws.Range("AF20").Value = 25
$ws.Range("AJ20").Value = 8
$ws.Range("AK20").Value = 19
$ws.Range("AM20").Value = 16
$ws.Range("AP20").Value = 22
$ws.Range("AQ20").Value = 13
$ws.Range("AS20").Value = 20
$ws.Range("AV20").Value = 8
$ws.Range("AW20").Value = 16
$ws.Range("AX20").Value = 23
$ws.Range("AZ20").Value = 22
$ws.Range("BA20").Value = 17
$ws.Range("BC20").Value = 25
$ws.Range("BF20").Value = "2016-06-12"

# Row 21
$ws.Range("D21").Value = 82
$ws.Range("E21").Value = 32
$ws.Range("G21").Value = 0.39
$ws.Range("I21").Value = 36.9
$ws.Range("J21").Value = 84
$ws.Range("O21").Value = 17.2
$ws.Range("P21").Value = 21.4
$ws.Range("Q21").Value = 0.805
$ws.Range("R21").Value = 10.4
$ws.Range("W21").Value = 5.7
$ws.Range("Z21").Value = 19.7
$ws.Range("AA21").Value = 18.5
$ws.Range("AB21").Value = 98.40000000000001
$ws.Range("AC21").Value = -2.7
$ws.Range("AI21").Value = 25
$ws.Range("AK21").Value = 26
$ws.Range("AO21").Value = 17
$ws.Range("AR21").Value = 16
$ws.Range("AV21").Value = 7
$ws.Range("BF21").Value = "2016-06-12"

# Row 22
$ws.Range("D22").Value = 82
$ws.Range("E22").Value = 55
$ws.Range("G22").Value = 0.671
$ws.Range("M22").Value = 23.7
$ws.Range("N22").Value = 0.349
$ws.Range("Q22").Value = 0.782
$ws.Range("T22").Value = 48.6
$ws.Range("X22").Value = 5.9
$ws.Range("AJ22").Value = 7
$ws.Range("AM22").Value = 17
$ws.Range("AN22").Value = 17
$ws.Range("AY22").Value = 11
$ws.Range("AZ22").Value = 16
$ws.Range("BF22").Value = "2016-06-12"

# Row 23
$ws.Range("D23").Value = 82
$ws.Range("F23").Value = 47
$ws.Range("G23").Value = 0.427
$ws.Range("H23").Value = 48.5
$ws.Range("I23").Value = 39.5
$ws.Range("J23").Value = 86.8
$ws.Range("K23").Value = 0.455
$ws.Range("M23").Value = 22.2
$ws.Range("N23").Value = 0.35
$ws.Range("O23").Value = 15.2
$ws.Range("P23").Value = 20.1
$ws.Range("Q23").Value = 0.757
$ws.Range("R23").Value = 10.3
$ws.Range("T23").Value = 43.3
$ws.Range("Y23").Value = 5.5
$ws.Range("AB23").Value = 102.1
$ws.Range("AH23").Value = 3
$ws.Range("AL23").Value = 23
$ws.Range("AN23").Value = 16
$ws.Range("AR23").Value = 18
$ws.Range("AT23").Value = 17
$ws.Range("AW23").Value = 12
$ws.Range("AY23").Value = 21
$ws.Range("AZ23").Value = 19
$ws.Range("BB23").Value = 18
$ws.Range("BC23").Value = 19
$ws.Range("BF23").Value = "2016-06-12"

# Row 24
$ws.Range("D24").Value = 82
$ws.Range("F24").Value = 72
$ws.Range("G24").Value = 0.122
$ws.Range("K24").Value = 0.431
$ws.Range("M24").Value = 27.5
$ws.Range("O24").Value = 15.7
$ws.Range("Q24").Value = 0.694
$ws.Range("R24").Value = 9.5
$ws.Range("U24").Value = 21.5
$ws.Range("V24").Value = 16.4
$ws.Range("AA24").Value = 19.2
$ws.Range("AC24").Value = -10.2
$ws.Range("AJ24").Value = 20
$ws.Range("AN24").Value = 24
$ws.Range("AQ24").Value = 28
$ws.Range("BA24").Value = 25
$ws.Range("BF24").Value = "2016-06-12"

# Row 25
$ws.Range("D25").Value = 82
$ws.Range("F25").Value = 59
$ws.Range("G25").Value = 0.28
$ws.Range("M25").Value = 25.8
$ws.Range("N25").Value = 0.348
$ws.Range("O25").Value = 17.5
$ws.Range("P25").Value = 23.2
$ws.Range("Q25").Value = 0.751
$ws.Range("U25").Value = 20.7
$ws.Range("W25").Value = 7.7
$ws.Range("Z25").Value = 22.7
$ws.Range("AA25").Value = 21.6
$ws.Range("AC25").Value = -6.7
$ws.Range("AL25").Value = 11
$ws.Range("AN25").Value = 18
$ws.Range("AO25").Value = 13
$ws.Range("AP25").Value = 13
$ws.Range("AQ25").Value = 20
$ws.Range("AR25").Value = 5
$ws.Range("AS25").Value = 15
$ws.Range("AU25").Value = 24
$ws.Range("AW25").Value = 17
$ws.Range("AY25").Value = 23
$ws.Range("BA25").Value = 6
$ws.Range("BF25").Value = "2016-06-12"

# Row 26
$ws.Range("D26").Value = 82
$ws.Range("E26").Value = 44
$ws.Range("G26").Value = 0.537
$ws.Range("J26").Value = 85.90000000000001
$ws.Range("K26").Value = 0.45
$ws.Range("M26").Value = 28.5
$ws.Range("P26").Value = 23
$ws.Range("Q26").Value = 0.754
$ws.Range("S26").Value = 33.9
$ws.Range("T26").Value = 45.5
$ws.Range("V26").Value = 14.6
$ws.Range("AC26").Value = 0.8
$ws.Range("AI26").Value = 10
$ws.Range("AK26").Value = 16
$ws.Range("AN26").Value = 4
$ws.Range("AO26").Value = 15
$ws.Range("AP26").Value = 14
$ws.Range("AQ26").Value = 19
$ws.Range("AR26").Value = 4
$ws.Range("AS26").Value = 8
$ws.Range("AV26").Value = 17
$ws.Range("BA26").Value = 22
$ws.Range("BF26").Value = "2016-06-12"

# Row 27
$ws.Range("D27").Value = 82
$ws.Range("F27").Value = 49
$ws.Range("G27").Value = 0.402
$ws.Range("J27").Value = 86.40000000000001
$ws.Range("N27").Value = 0.359
$ws.Range("O27").Value = 18.5
$ws.Range("P27").Value = 25.5
$ws.Range("Q27").Value = 0.725
$ws.Range("S27").Value = 33.7
$ws.Range("V27").Value = 16.2
$ws.Range("W27").Value = 8.9
$ws.Range("AC27").Value = -2.5
$ws.Range("AJ27").Value = 6
$ws.Range("AN27").Value = 10
$ws.Range("AO27").Value = 9
$ws.Range("AP27").Value = 6
$ws.Range("AT27").Value = 11
$ws.Range("AU27").Value = 4
$ws.Range("BA27").Value = 4
$ws.Range("BB27").Value = 3
$ws.Range("BF27").Value = "2016-06-12"

# Row 28
$ws.Range("D28").Value = 82
$ws.Range("F28").Value = 15
$ws.Range("G28").Value = 0.8169999999999999
$ws.Range("I28").Value = 40.1
$ws.Range("J28").Value = 82.90000000000001
$ws.Range("M28").Value = 18.5
$ws.Range("N28").Value = 0.375
$ws.Range("P28").Value = 20.4
$ws.Range("Q28").Value = 0.803
$ws.Range("S28").Value = 34.5
$ws.Range("T28").Value = 43.9
$ws.Range("W28").Value = 8.300000000000001
$ws.Range("AC28").Value = 10.6
$ws.Range("AR28").Value = 25
$ws.Range("AT28").Value = 15
$ws.Range("AV28").Value = 4
$ws.Range("BB28").Value = 10
$ws.Range("BC28").Value = 2
$ws.Range("BF28").Value = "2016-06-12"

# Row 29
$ws.Range("D29").Value = 82
$ws.Range("E29").Value = 56
$ws.Range("G29").Value = 0.6830000000000001
$ws.Range("L29").Value = 8.6
$ws.Range("O29").Value = 20.8
$ws.Range("P29").Value = 26.7
$ws.Range("Q29").Value = 0.777
$ws.Range("S29").Value = 33.2
$ws.Range("T29").Value = 43.4
$ws.Range("V29").Value = 13.1
$ws.Range("AA29").Value = 22
$ws.Range("AH29").Value = 23
$ws.Range("AK29").Value = 15
$ws.Range("AN29").Value = 5
$ws.Range("AR29").Value = 19
$ws.Range("AS29").Value = 18
$ws.Range("AV29").Value = 5
$ws.Range("AW29").Value = 15
$ws.Range("BA29").Value = 3
$ws.Range("BF29").Value = "2016-06-12"

# Row 30
$ws.Range("D30").Value = 82
$ws.Range("F30").Value = 42
$ws.Range("G30").Value = 0.488
$ws.Range("J30").Value = 80.40000000000001
$ws.Range("K30").Value = 0.449
$ws.Range("M30").Value = 23.9
$ws.Range("N30").Value = 0.355
$ws.Range("R30").Value = 10.7
$ws.Range("T30").Value = 43.2
$ws.Range("V30").Value = 14.9
$ws.Range("Z30").Value = 20.2
$ws.Range("AB30").Value = 97.7
$ws.Range("AC30").Value = 1.8
$ws.Range("AE30").Value = 19
$ws.Range("AF30").Value = 19
$ws.Range("AG30").Value = 19
$ws.Range("AK30").Value = 18
$ws.Range("AO30").Value = 19
$ws.Range("AQ30").Value = 24
$ws.Range("AT30").Value = 18
$ws.Range("AV30").Value = 21
$ws.Range("AW30").Value = 18
$ws.Range("BA30").Value = 18
$ws.Range("BF30").Value = "2016-06-12"

# Row 31
$ws.Range("D31").Value = 82
$ws.Range("E31").Value = 41
$ws.Range("G31").Value = 0.5
$ws.Range("I31").Value = 39.5
$ws.Range("K31").Value = 0.46
$ws.Range("L31").Value = 8.6
$ws.Range("M31").Value = 24.2
$ws.Range("N31").Value = 0.358
$ws.Range("Q31").Value = 0.73
$ws.Range("R31").Value = 9.1
$ws.Range("V31").Value = 14.5
$ws.Range("Z31").Value = 20.8
$ws.Range("AB31").Value = 104.1
$ws.Range("AF31").Value = 17
$ws.Range("AG31").Value = 17
$ws.Range("AI31").Value = 6
$ws.Range("AN31").Value = 11
$ws.Range("AU31").Value = 5
$ws.Range("AZ31").Value = 21
$ws.Range("BF31").Value = "2016-06-12"
